$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22: new journal entry ---
$ws.Range("B22").NumberFormat = $ws.Range("B21").NumberFormat
$ws.Range("B22").Value = 43632
$ws.Range("C22").Value = "Conception"
$ws.Range("D22").Value = "Adaptation du template trouvé au modèle MVC. Création de l'accès à la page d'accueil. Cet accès n'est actuellement plus possible tant que je n'aurai pas résolu un problème de chargement de la bonne page."
$ws.Range("E22").Value = 6
$ws.Rows.Item(22).RowHeight = 75

# --- Row 23: new journal entry ---
$ws.Range("B23").NumberFormat = $ws.Range("B21").NumberFormat
$ws.Range("B23").Value = 43601
$ws.Range("C23").Value = "Conception"
$ws.Range("D23").Value = "Modification de mes entrées dans la base de données afin que celle-ci soient en anglais."
$ws.Range("E23").Value = 0.75
$ws.Rows.Item(23).RowHeight = 30

# --- Update the view to reflect where the author was working ---
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("E23").Select()
